$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the newly-added "Nugg" values for the first four package rows ---
$ws.Range("E2").Value = "HCVN_ZABAP_UTIL"
$ws.Range("E3").Value = "HCVN_ZABAPGIT"
$ws.Range("E4").Value = "HCVN_ZBC_ACTIVITY_LOG"
$ws.Range("E5").Value = "HCVN_ZBC_GENERAL"

# --- 2. Highlight rows 2-5 (A:E) with the new blue accent fill ---
$ws.Range("A2:E5").Interior.Color = 15652797

# --- 3. Re-color the "Nugg" column for rows 31-50 and 53 to match the yellow
#         highlight already used elsewhere in the sheet ---
$ws.Range("E31:E50").Interior.Color = 65535
$ws.Range("E53").Interior.Color = 65535

# --- 4. Recolor row 51 entirely using the same gray fill used on rows 13/15,
#         copying the format so the existing style/fill gets reused verbatim ---
$ws.Range("A13:E13").Copy() | Out-Null
$ws.Range("A51:E51").PasteSpecial(-4122) | Out-Null

# --- 5. Update the active selection / scroll position ---
$ws.Range("G11").Select() | Out-Null
